$wb = $excel.ActiveWorkbook

# This script applies the numeric corrections described in the commit diff
# ("chore: update Sheets via scheduled runner") across the eight Leve-profit
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Each worksheet has an
# identical A1:N141 layout; only specific H/I/J/K/L/M/N cells on specific rows
# change value, and a couple of cells are added or cleared entirely.

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H15").Value = 1139.3572
$ws.Range("I15").Value = 1139.3572
$ws.Range("K15").Value = 3418.0716
$ws.Range("M15").Value = -3249.0716
$ws.Range("H113").Value = 3828.4285
$ws.Range("I113").Value = 3266.6667
$ws.Range("K113").Value = 3266.6667
$ws.Range("M113").Value = -12.66670000000022
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 4673.8394
$ws.Range("I132").Value = 5627.622
$ws.Range("K132").Value = 16882.866
$ws.Range("M132").Value = -14352.866

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H5").Value = 2968.6667
$ws.Range("I5").Value = 44.2
$ws.Range("J5").Value = 6624.25
$ws.Range("K5").Value = 44.2
$ws.Range("L5").Value = 6624.25
$ws.Range("M5").Value = 67.8
$ws.Range("N5").Value = -6848.25
$ws.Range("H45").Value = 3039.3333
$ws.Range("J45").Value = 4499.75
$ws.Range("L45").Value = 4499.75
$ws.Range("N45").Value = -5253.75
$ws.Range("H46").Value = 7562.5
$ws.Range("I46").Value = 4625
$ws.Range("K46").Value = 4625
$ws.Range("M46").Value = -4306
$ws.Range("H110").Value = 872.5333000000001
$ws.Range("I110").Value = 884.8570999999999
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 884.8570999999999
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 1160.1429
$ws.Range("N110").Value = -4790
$ws.Range("H139").Value = 89718.086
$ws.Range("J139").Value = 89718.086
$ws.Range("L139").Value = 89718.086
$ws.Range("N139").Value = -99998.086

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H4").Value = 2968.6667
$ws.Range("I4").Value = 44.2
$ws.Range("J4").Value = 6624.25
$ws.Range("K4").Value = 44.2
$ws.Range("L4").Value = 6624.25
$ws.Range("M4").Value = 70.8
$ws.Range("N4").Value = -6854.25
$ws.Range("H97").Value = 19920.223
$ws.Range("I97").Value = 14910.25
$ws.Range("K97").Value = 14910.25
$ws.Range("M97").Value = -13919.25
$ws.Range("H107").Value = 990.1429000000001
$ws.Range("I107").Value = 974.4167
$ws.Range("J107").Value = 1084.5
$ws.Range("K107").Value = 974.4167
$ws.Range("L107").Value = 1084.5
$ws.Range("M107").Value = 945.5833
$ws.Range("N107").Value = -4924.5

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 1043936.25
$ws.Range("I31").Value = 1545572.2
$ws.Range("K31").Value = 1545572.2
$ws.Range("M31").Value = -1545277.2
$ws.Range("H34").Value = 1043936.25
$ws.Range("I34").Value = 1545572.2
$ws.Range("K34").Value = 1545572.2
$ws.Range("M34").Value = -1545370.2
$ws.Range("H68").Value = 42000
$ws.Range("J68").Value = 42000
$ws.Range("L68").Value = 42000
$ws.Range("N68").Value = -43498
$ws.Range("H71").Value = 42000
$ws.Range("J71").Value = 42000
$ws.Range("L71").Value = 126000
$ws.Range("N71").Value = -133488
$ws.Range("H74").Value = 45000
$ws.Range("J74").Value = 45000
$ws.Range("L74").Value = 45000
$ws.Range("N74").Value = -46748
$ws.Range("H77").Value = 45000
$ws.Range("J77").Value = 45000
$ws.Range("L77").Value = 135000
$ws.Range("N77").Value = -143736
$ws.Range("H105").Value = 10218.909
$ws.Range("J105").Value = 1794.6666
$ws.Range("L105").Value = 1794.6666
$ws.Range("N105").Value = -5288.6666
$ws.Range("H107").Value = 961.5714
$ws.Range("I107").Value = 768.05554
$ws.Range("J107").Value = 2122.6667
$ws.Range("K107").Value = 768.05554
$ws.Range("L107").Value = 2122.6667
$ws.Range("M107").Value = 1151.94446
$ws.Range("N107").Value = -5962.6667
$ws.Range("H110").Value = 48851
$ws.Range("J110").Value = 48851
$ws.Range("L110").Value = 48851
$ws.Range("N110").Value = -57031
$ws.Range("H115").Value = 12145
$ws.Range("J115").Value = 12145
$ws.Range("L115").Value = 12145
$ws.Range("N115").Value = -14495

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H38").Value = 172.20833
$ws.Range("J38").Value = 126.46667
$ws.Range("L38").Value = 379.40001
$ws.Range("N38").Value = -1073.40001
$ws.Range("H39").Value = 600
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H62").Value = 2156.3333
$ws.Range("J62").Value = 2984.5
$ws.Range("L62").Value = 8953.5
$ws.Range("N62").Value = -10325.5
$ws.Range("H65").Value = 2156.3333
$ws.Range("J65").Value = 2984.5
$ws.Range("L65").Value = 26860.5
$ws.Range("N65").Value = -33724.5

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H97").Value = 1223.7084
$ws.Range("I97").Value = 1169.7333
$ws.Range("J97").Value = 1313.6666
$ws.Range("K97").Value = 1169.7333
$ws.Range("L97").Value = 1313.6666
$ws.Range("M97").Value = -673.7333000000001
$ws.Range("N97").Value = -2305.6666
$ws.Range("H132").Value = 24312.111
$ws.Range("I132").Value = 25564.941
$ws.Range("K132").Value = 76694.823
$ws.Range("M132").Value = -74164.823

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H33").Value = 8959.333000000001
$ws.Range("J33").Value = 8959.333000000001
$ws.Range("L33").Value = 8959.333000000001
$ws.Range("N33").Value = -9539.333000000001
$ws.Range("H99").Value = 43723.5
$ws.Range("I99").Value = 43298
$ws.Range("K99").Value = 43298
$ws.Range("M99").Value = -40303
$ws.Range("H132").Value = 1669088.1
$ws.Range("I132").Value = 2022076.5
$ws.Range("J132").Value = 4999.7144
$ws.Range("K132").Value = 6066229.5
$ws.Range("L132").Value = 14999.1432
$ws.Range("M132").Value = -6063699.5
$ws.Range("N132").Value = -20059.1432

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H107").Value = 38464784
$ws.Range("J107").Value = 4630.6924
$ws.Range("L107").Value = 13892.0772
$ws.Range("N107").Value = -17732.0772
$ws.Range("H126").Value = 2032.4286
$ws.Range("I126").Value = 1763.4736
$ws.Range("J126").Value = 2600.2222
$ws.Range("K126").Value = 5290.4208
$ws.Range("L126").Value = 7800.6666
$ws.Range("M126").Value = -2820.4208
$ws.Range("N126").Value = -12740.6666
$ws.Range("H132").Value = 9805022
$ws.Range("I132").Value = 11905876
$ws.Range("J132").Value = 1033.3334
$ws.Range("K132").Value = 35717628
$ws.Range("L132").Value = 3100.0002
$ws.Range("M132").Value = -35715098
$ws.Range("N132").Value = -8160.0002
